$d = $word.ActiveDocument

# Helper: rewrite the paragraph currently at position $idx with a single
# clean run containing $text, discarding any proofErr spell/grammar-check
# markers and multi-run splits that were present in the original paragraph.
function Clean-ParagraphAt {
    param($idx, $text)

    $para = $d.Paragraphs.Item($idx)
    # Insert a brand-new (empty, markup-free) paragraph immediately before
    # the existing one; the new paragraph takes index $idx and the old
    # (possibly "dirty") paragraph shifts to $idx + 1.
    $null = $para.Range.InsertParagraphBefore()
    $d.Paragraphs.Item($idx).Range.Text = $text
    # Remove the old paragraph (and its end-of-paragraph mark) entirely so
    # none of its runs/proofErr markers survive.
    $d.Paragraphs.Item($idx + 1).Range.Delete()
}

# Paragraph 1: "Dat_use ..." – merge the two split runs (and drop the
# spellStart/spellEnd proofErr pair) into a single clean run.
Clean-ParagraphAt 1 "Dat_use – this is the raw cleaned data, not scaled or centered. At the commune scale"

# Paragraph 2: "Dat1 ..." – merge the two split runs (and drop the
# gramStart/gramEnd proofErr pair) into a single clean run.
Clean-ParagraphAt 2 "Dat1 – this is the above data but scaled"

# Insert the new "Dat_prov" paragraph right after "Dat1" and before "Dat2".
$p2 = $d.Paragraphs.Item(2)
$null = $p2.Range.InsertParagraphAfter()
$d.Paragraphs.Item(3).Range.Text = "Dat_prov – this is the unscaled data aggregated up to the province level"

# Paragraph (now 4): "Dat2 ..." – merge the two split runs (and drop the
# gramStart/gramEnd proofErr pair) into a single clean run.
Clean-ParagraphAt 4 "Dat2 – this is the scaled data aggregated up to the Province level"

# Paragraph (now 5): "Dat_cat ..." – merge the two split runs (and drop the
# spellStart/spellEnd proofErr pair) into a single clean run.
Clean-ParagraphAt 5 "Dat_cat – this is the provincial level data but with some (most) of the variables transformed to categorical variables. This was done by splitting the data into those observations on the left of the mean (“low”) and those on the right of the mean (“high”)"
